$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117, shifting existing rows 117-119 down to 118-120
$ws.Rows.Item(117).Insert()

# Populate new row 117 with data matching the surrounding dataset pattern,
# with the updated date/volume/price values for this record
$ws.Cells.Item(117, 1).Value = 10
$ws.Cells.Item(117, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(117, 3).Value = "La Araucanía"
$ws.Cells.Item(117, 4).Value = 44448
$ws.Cells.Item(117, 5).Value = 9
$ws.Cells.Item(117, 6).Value = 100112005
$ws.Cells.Item(117, 7).Value = "Puerro"
$ws.Cells.Item(117, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 60
$ws.Cells.Item(117, 11).Value = 8000
$ws.Cells.Item(117, 12).Value = 8000
$ws.Cells.Item(117, 13).Value = 8000
$ws.Cells.Item(117, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(117, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(117, 16).Value = 667
$ws.Cells.Item(117, 17).Value = 12
$ws.Cells.Item(117, 18).Value = "Hortaliza"
